$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '67.617.09'
$ws.Range("E2").Value = '  +0.16%  '

# Row 3
$ws.Range("D3").Value = '3.304.10'
$ws.Range("E3").Value = '  -2.08%  '

# Row 4
$ws.Range("E4").Value = '  -0.01%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.81'
$ws.Range("E5").Value = '  -1.79%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '174.32'
$ws.Range("E6").Value = '  -6.86%  '

# Row 7
$ws.Range("E7").Value = '  +0.04%  '

# Row 8
$ws.Range("E8").Value = '  -1.85%  '

# Row 9
$ws.Range("D9").Value = '3.302.76'
$ws.Range("E9").Value = '  -1.90%  '

# Row 10
$ws.Range("E10").Value = '  -4.70%  '

# Row 11
$ws.Range("E11").Value = '  -2.22%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '45.22'
$ws.Range("E12").Value = '  -4.67%  '

# Row 13
$ws.Range("E13").Value = '  -2.27%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '667.84'
$ws.Range("E14").Value = '  +5.23%  '

# Row 15
$ws.Range("D15").Value = '3.840.35'
$ws.Range("E15").Value = '  -2.02%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.35'
$ws.Range("E16").Value = '  -3.00%  '

# Row 17
$ws.Range("D17").Value = '67.656.07'
$ws.Range("E17").Value = '  +0.10%  '

# Row 18
$ws.Range("E18").Value = '  -0.83%  '

# Row 19
$ws.Range("D19").Value = '3.304.28'
$ws.Range("E19").Value = '  -2.18%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.37'
$ws.Range("E20").Value = '  -3.50%  '

# Row 21
$ws.Range("E21").Value = '  -3.09%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.884'
$ws.Range("E22").Value = '  -2.80%  '

# Row 23
$ws.Range("B23").Value = 'Toncoin'
$ws.Range("C23").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.37'
$ws.Range("E23").Value = '  +5.13%  '

# Row 24
$ws.Range("B24").Value = 'InternetComputer(DFINITY)'
$ws.Range("C24").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '17.06'
$ws.Range("E24").Value = '  -5.15%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '97.12'
$ws.Range("E25").Value = '  -2.23%  '

# Row 26
$ws.Range("E26").Value = '  -4.44%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.67'
$ws.Range("E27").Value = '  -6.55%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.16'
$ws.Range("E28").Value = '  -5.51%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '32.98'
$ws.Range("E29").Value = '  +1.09%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.36'
$ws.Range("E30").Value = '  -3.71%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.07'
$ws.Range("E31").Value = '  +1.65%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '590.59'
$ws.Range("E32").Value = '  -2.98%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '10.91'
$ws.Range("E33").Value = '  -1.58%  '

# Row 34
$ws.Range("E34").Value = '  -2.54%  '

# Row 35
$ws.Range("D35").Value = '3.722.58'
$ws.Range("E35").Value = '  -6.78%  '

# Row 36
$ws.Range("E36").Value = '  +0.11%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.32'
$ws.Range("E37").Value = '  -12.46%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '55.21'
$ws.Range("E38").Value = '  -1.40%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.132'
$ws.Range("E39").Value = '  +0.21%  '

# Row 40
$ws.Range("B40").Value = 'Fetch.AI'
$ws.Range("C40").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.61'
$ws.Range("E40").Value = '  -8.15%  '

# Row 41
$ws.Range("B41").Value = 'InjectiveProtocol'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '32.12'
$ws.Range("E41").Value = '  -4.57%  '

# Row 42
$ws.Range("E42").Value = '  -4.61%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.28'
$ws.Range("E43").Value = '  -3.20%  '

# Row 44
$ws.Range("E44").Value = '  -5.64%  '

# Row 45
$ws.Range("E45").Value = '  -3.81%  '

# Row 46
$ws.Range("E46").Value = '  -4.28%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.59'
$ws.Range("E47").Value = '  +0.40%  '

# Row 48
$ws.Range("E48").Value = '  -2.27%  '

# Row 49
$ws.Range("E49").Value = '  +0.17%  '

# Row 50
$ws.Range("E50").Value = '  -2.34%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '129.36'
$ws.Range("E51").Value = '  +0.95%  '
